# "opex without building rent"
# Update the per-technology OPEX-per-year figures on the OPEX sheet
# (B3:B15). B2 (the average) is a formula (=AVERAGE(B3:B15)) and will
# recalc automatically, as will every downstream consumer: the
# Revenue sheet's U/V/W/X/Y/Z/AA/.../AH columns (which all key off
# OPEX!$B$15) and the "Approx. OPEX per year" bar chart (chart2, whose
# series is the live range OPEX!$B$2:$B$15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OPEX")

$ws.Range("B3").Value = 19283.042118970367   # FTTC_GPON_25
$ws.Range("B4").Value = 43709.823692218881   # FTTB_XGPON_50
$ws.Range("B5").Value = 44248.4220525376     # FTTB_UDWDM_50
$ws.Range("B6").Value = 10540.776519765759   # FTTH_UDWDM_100
$ws.Range("B7").Value = 11938.918689039821   # FTTH_XGPON_100
$ws.Range("B8").Value = 50983.71612182937    # FTTC_GPON_100
$ws.Range("B9").Value = 47726.863400003029   # FTTB_XGPON_100
$ws.Range("B10").Value = 48220.289169209609  # FTTB_UDWDM_100
$ws.Range("B11").Value = 37682.867713520049  # FTTC_Hybridpon_25
$ws.Range("B12").Value = 44827.8721002112    # FTTB_Hybridpon_50
$ws.Range("B13").Value = 12139.814930199653  # FTTH_Hybridpon_100 (unchanged)
$ws.Range("B14").Value = 47419.574760442287  # FTTC_Hybridpon_100
$ws.Range("B15").Value = 50088.476920599649  # FTTB_Hybridpon_100

# Restore the cursor position on the OPEX sheet the way it was left
# after the edit (selection moved from K3 to B38).
$ws.Range("B38").Select()

$excel.CalculateFull()
